# Updates price (D) and 1h volume % (E) columns for the crypto symbol list,
# mirroring the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Leading apostrophe forces Excel to store the value as literal text,
# matching the original inlineStr cell type (instead of being auto-parsed
# into a number / percentage).
$updates = @{
    "D2" = "'301.19"
    "E2" = "'-0.87%"
    "D3" = "'31.37"
    "E3" = "'-2.09%"
    "D4" = "'5.112"
    "E4" = "'-2.30%"
    "D5" = "'0.07359"
    "E5" = "'-2.60%"
    "D6" = "'2.182"
    "E6" = "'42.33%"
    "D7" = "'7.929"
    "E7" = "'-0.30%"
    "D8" = "'3.820"
    "E8" = "'-0.70%"
    "D9" = "'0.9181"
    "E9" = "'-1.13%"
    "D10" = "'0.1705"
    "E10" = "'0.56%"
    "D11" = "'0.07487"
    "E11" = "'-4.89%"
    "D12" = "'0.08131"
    "E12" = "'1.27%"
    "D13" = "'0.03022"
    "E13" = "'-0.17%"
    "D14" = "'0.09916"
    "E14" = "'0.03%"
    "D15" = "'0.001496"
    "E15" = "'-1.44%"
    "D16" = "'0.006157"
    "E16" = "'-3.31%"
    "D17" = "'3.458"
    "E17" = "'0.05%"
    "D18" = "'2.224"
    "E18" = "'-0.45%"
    "D19" = "'0.3281"
    "E19" = "'-0.54%"
    "D20" = "'0.1350"
    "E20" = "'1.30%"
    "D21" = "'4.653"
    "E21" = "'4.55%"
    "D22" = "'0.04636"
    "E22" = "'0.96%"
    "D23" = "'0.1567"
    "E23" = "'-3.15%"
    "D24" = "'0.001226"
    "E24" = "'0.71%"
    "D25" = "'0.004475"
    "E25" = "'-0.18%"
    "D26" = "'0.0001299"
    "E26" = "'-6.99%"
    "D27" = "'0.0003427"
    "E27" = "'92.41%"
    "E39" = "'1.26%"
    "D40" = "'0.04501"
    "E40" = "'0.08%"
    "D41" = "'0.007234"
    "E41" = "'3.57%"
    "D42" = "'0.1343"
    "E42" = "'-1.23%"
    "D43" = "'0.002228"
    "E43" = "'7.37%"
    "D44" = "'0.01062"
    "E44" = "'-22.52%"
    "D45" = "'0.00006303"
    "E45" = "'2.46%"
    "D46" = "'0.8212"
    "E46" = "'14.20%"
    "E47" = "'-22.94%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
